$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 8.98775
$ws.Range("C2").Value = 4.81645
$ws.Range("D2").Value = 34.36245
$ws.Range("F2").Value = 27.96965
$ws.Range("G2").Value = 27.96965
$ws.Range("K2").Value = 622.0352
$ws.Range("L2").Value = 27.97
$ws.Range("M2").Value = 594.0652
$ws.Range("N2").Value = 9.794799999999999
$ws.Range("O2").Value = 584.2704
$ws.Range("F3").Value = 34.21
$ws.Range("G3").Value = 32.538
$ws.Range("K3").Value = 40.0614
$ws.Range("L3").Value = 32.538
$ws.Range("B4").Value = 40.966
$ws.Range("C4").Value = 11.645
$ws.Range("D4").Value = 40.083
$ws.Range("F4").Value = 39.014
$ws.Range("G4").Value = 39.01360821917808
$ws.Range("K4").Value = 39.014
$ws.Range("L4").Value = 39.014

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 8.982899999999999
$ws.Range("C2").Value = 4.54085
$ws.Range("D2").Value = 33.99040000000001
$ws.Range("E2").Value = 0.7522
$ws.Range("F2").Value = 27.94
$ws.Range("G2").Value = 27.94
$ws.Range("K2").Value = 391.9772
$ws.Range("L2").Value = 27.94
$ws.Range("M2").Value = 364.0371999999999
$ws.Range("N2").Value = 9.5928
$ws.Range("O2").Value = 354.4444
$ws.Range("F3").Value = 34.115
$ws.Range("G3").Value = 32.074
$ws.Range("K3").Value = 84.0728
$ws.Range("L3").Value = 32.074
$ws.Range("B4").Value = 32.898
$ws.Range("C4").Value = 6.722
$ws.Range("D4").Value = 50.187
$ws.Range("E4").Value = 1.661
$ws.Range("F4").Value = 36.111
$ws.Range("G4").Value = 36.11118082191781
$ws.Range("K4").Value = 146.3526
$ws.Range("L4").Value = 36.111
$ws.Range("M4").Value = 110.2416
$ws.Range("N4").Value = 1.0082
$ws.Range("O4").Value = 109.2334

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 10.85825
$ws.Range("C2").Value = 4.805750000000001
$ws.Range("D2").Value = 32.34480000000001
$ws.Range("E2").Value = 0.9064499999999999
$ws.Range("F2").Value = 28.53235
$ws.Range("G2").Value = 28.53235
$ws.Range("K2").Value = 77.23400000000001
$ws.Range("L2").Value = 28.532
$ws.Range("M2").Value = 48.702
$ws.Range("N2").Value = 7.4026
$ws.Range("O2").Value = 41.29940000000001
$ws.Range("B3").Value = 24.126
$ws.Range("C3").Value = 14.389
$ws.Range("D3").Value = 35.425
$ws.Range("F3").Value = 36.084
$ws.Range("G3").Value = 33.801
$ws.Range("K3").Value = 33.801
$ws.Range("L3").Value = 33.801
$ws.Range("B4").Value = 40.966
$ws.Range("C4").Value = 11.645
$ws.Range("D4").Value = 39.838
$ws.Range("F4").Value = 39.009
$ws.Range("G4").Value = 39.00857260273973
$ws.Range("K4").Value = 39.02999999999999
$ws.Range("L4").Value = 39.009

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 11.7456
$ws.Range("C2").Value = 4.97875
$ws.Range("D2").Value = 32.42345
$ws.Range("E2").Value = 0.9064499999999999
$ws.Range("F2").Value = 28.83395
$ws.Range("G2").Value = 28.83395
$ws.Range("K2").Value = 30.1378
$ws.Range("L2").Value = 28.834
$ws.Range("M2").Value = 1.3038
$ws.Range("N2").Value = 1.3038
$ws.Range("B3").Value = 24.126
$ws.Range("C3").Value = 14.389
$ws.Range("D3").Value = 35.425
$ws.Range("F3").Value = 36.084
$ws.Range("G3").Value = 33.801
$ws.Range("K3").Value = 33.801
$ws.Range("L3").Value = 33.801
$ws.Range("B4").Value = 40.966
$ws.Range("C4").Value = 11.645
$ws.Range("D4").Value = 40.083
$ws.Range("F4").Value = 39.014
$ws.Range("G4").Value = 39.01360821917808
$ws.Range("K4").Value = 39.014
$ws.Range("L4").Value = 39.014

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 11.4189
$ws.Range("C2").Value = 4.8712
$ws.Range("D2").Value = 32.6733
$ws.Range("E2").Value = 0.7324999999999999
$ws.Range("F2").Value = 28.72045
$ws.Range("G2").Value = 28.72045
$ws.Range("K2").Value = 70.23179999999999
$ws.Range("L2").Value = 28.72
$ws.Range("M2").Value = 41.5118
$ws.Range("N2").Value = 2.3116
$ws.Range("O2").Value = 39.2002
$ws.Range("F3").Value = 36.058
$ws.Range("G3").Value = 32.505
$ws.Range("K3").Value = 32.6428
$ws.Range("L3").Value = 32.505
$ws.Range("B4").Value = 40.966
$ws.Range("C4").Value = 11.645
$ws.Range("D4").Value = 40.083
$ws.Range("F4").Value = 39.014
$ws.Range("G4").Value = 39.01360821917808
$ws.Range("K4").Value = 39.014
$ws.Range("L4").Value = 39.014
